$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters to update for each data row (G stays 1, H stays blank - unchanged)
$cols = @("B", "C", "D", "E", "F", "I", "J", "K", "L", "M", "N")

# Map: row number -> array of new values aligned with $cols
$rowValues = @{}
$rowValues[2] = @(1.02, 1.058827272464505, 1.059959359295507, 1.063457380140247, 1.072421245231305, 1.044036366699587, 1.063816471772534, 1.062687442571183, 1.066175967479244, 1.075115804479646, 1.065327213183905)
$rowValues[3] = @(1.02, 1.060443490833719, 1.061190878046065, 1.064781522748067, 1.073831823864521, 1.044442690264865, 1.065082511300405, 1.063732238804386, 1.067313847316746, 1.076341661068609, 1.066595050633123)
$rowValues[4] = @(1.02, 1.061487391272947, 1.061985926517937, 1.065637419781649, 1.07474331930633, 1.044703118963766, 1.06589945575311, 1.064405904229036, 1.068048674804747, 1.077133099185782, 1.067413155240656)
$rowValues[5] = @(1.02, 1.061925801570217, 1.06231973387656, 1.065997026643344, 1.075126221195211, 1.0448120098, 1.066242363570668, 1.064688547430348, 1.068357253121175, 1.077465401614874, 1.067756550026647)
$rowValues[6] = @(1.02, 1.061999386745809, 1.06237575650648, 1.066057393945775, 1.075190495252838, 1.0448302583313, 1.066299908103057, 1.064735971461321, 1.068409044817754, 1.077521172316916, 1.067814176278872)
$rowValues[7] = @(1.02, 1.061493251069679, 1.061990388553001, 1.065642225691618, 1.074748436793098, 1.044704576296258, 1.065904039798615, 1.064409683136074, 1.068052799383775, 1.07713754106127, 1.067417745796031)
$rowValues[8] = @(1.02, 1.05937388044584, 1.060375938790504, 1.063905072764069, 1.07289821844183, 1.044174202803833, 1.064244809714697, 1.06304103425505, 1.06656082433905, 1.075530459395647, 1.065756159415078)
$rowValues[9] = @(1.02, 1.055624286389638, 1.057516801455672, 1.060836732914596, 1.069628067381169, 1.043220426891146, 1.06130333800961, 1.060610750091685, 1.063920363278827, 1.072684712005929, 1.062810510482778)
$rowValues[10] = @(1.02, 1.053113868460232, 1.055600718267628, 1.058785913796267, 1.06744092407737, 1.042571502808724, 1.05932999776765, 1.058977717100004, 1.06215202876633, 1.070777825041792, 1.060834367871329)
$rowValues[11] = @(1.02, 1.052024158438887, 1.054768574005567, 1.057896553413974, 1.066492099198257, 1.042287371983355, 1.05847248483343, 1.0582674659612, 1.061384336732891, 1.069949726087418, 1.059975637170385)
$rowValues[12] = @(1.02, 1.051618976786898, 1.05445910005913, 1.057565996263662, 1.066139387886373, 1.04218135761783, 1.058153500077215, 1.058003168458441, 1.061098875576436, 1.069641764204074, 1.059656199419228)
$rowValues[13] = @(1.02, 1.051705908589462, 1.054525500525662, 1.057636911477425, 1.066215058264107, 1.042204119632853, 1.058221944635508, 1.058059882945755, 1.061160121929843, 1.069707839900841, 1.059724741176644)
$rowValues[14] = @(1.02, 1.051990674528285, 1.054743000559548, 1.05786923376588, 1.066462949640692, 1.042278618523293, 1.058446127018175, 1.05824562885304, 1.06136074671051, 1.069924277436398, 1.059949241924008)
$rowValues[15] = @(1.02, 1.052166072897272, 1.054876959117495, 1.058012347262635, 1.066615646937379, 1.042324456633908, 1.058584191126693, 1.058360009397539, 1.06148431748205, 1.070057582605549, 1.060087502099396)
$rowValues[16] = @(1.02, 1.053186131402214, 1.055655892289247, 1.05884490880255, 1.067503856290663, 1.042590293138208, 1.05938684318514, 1.059024787413201, 1.062202935407318, 1.070832731918441, 1.060891294015832)
$rowValues[17] = @(1.02, 1.053825260663578, 1.056143830368063, 1.059366788076707, 1.068060524735172, 1.0427562016414, 1.059889504427697, 1.059440940085713, 1.062653167371182, 1.071318313387269, 1.06139466909501)
$rowValues[18] = @(1.02, 1.054197795930849, 1.056428198837311, 1.059671062435565, 1.068385049104633, 1.042852670299658, 1.060182405177214, 1.059683372713301, 1.062915588378592, 1.071601313394813, 1.061687985797188)
$rowValues[19] = @(1.02, 1.054324777421271, 1.056525121077027, 1.059774790473906, 1.068495674792695, 1.042885512333658, 1.060282227328407, 1.059765984937809, 1.063005034857478, 1.07169777006229, 1.061787949707285)
$rowValues[20] = @(1.02, 1.053756714915977, 1.056091503858018, 1.059310808765936, 1.068000817219351, 1.042738432590926, 1.059835604010127, 1.059396322151666, 1.062604881619107, 1.071266239087668, 1.061340692132664)
$rowValues[21] = @(1.02, 1.051906829603681, 1.054678962732431, 1.057800826474974, 1.06638995948253, 1.042256693611335, 1.058380123838206, 1.05819094458374, 1.0613016762205, 1.06986055222793, 1.059883145011953)
$rowValues[22] = @(1.02, 1.050741327274853, 1.053788648172895, 1.056850227115983, 1.065375551047981, 1.041951052230498, 1.05746230346058, 1.057430303135253, 1.06048052446532, 1.068974601902286, 1.058964021224116)
$rowValues[23] = @(1.02, 1.051359414144945, 1.054260831340501, 1.057354275360182, 1.065913462494656, 1.042113340607369, 1.057949116576282, 1.057833798834938, 1.060916003113007, 1.069444466339718, 1.059451525670279)
$rowValues[24] = @(1.02, 1.05378768858617, 1.056115148689768, 1.059336103805682, 1.06802779700401, 1.0427464625962, 1.059859960194353, 1.059416484009388, 1.0626267004645, 1.071289769937489, 1.061365082905466)
$rowValues[25] = @(1.02, 1.056595485338559, 1.058257688760585, 1.061630871114318, 1.070474689115292, 1.043469292454112, 1.062065923277499, 1.061241273111823, 1.064604374361861, 1.073422089917419, 1.063574178709217)

foreach ($r in $rowValues.Keys) {
    $vals = $rowValues[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}

Write-Output "Updated $($rowValues.Keys.Count) rows"
